# Add a new worksheet "04 04 2022" after the existing "21 03 2022" sheet,
# cloning its layout (headers, "Midi"/"Soir" blocks, day numbering) and then
# tweaking the 04/04 menu: the "Entrée" for day 1 of the midi service becomes
# "Crêpes jambon" instead of a plain number, and the evening ("Soir") block
# has no Friday (F) column of data like the original sheet did.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "04 04 2022"

# ---- Row 1 : header row -------------------------------------------------
$ws2.Range("A1").Value = "Jour"
$ws2.Range("A1").NumberFormat = "dd/mm/yy"
$ws2.Range("B1").Value = "Lundi"
$ws2.Range("B1").NumberFormat = "dd/mm/yy"
$ws2.Range("C1").Value = "Mardi"
$ws2.Range("C1").NumberFormat = "dd/mm/yy"
$ws2.Range("D1").Value = "Mercredi"
$ws2.Range("D1").NumberFormat = "dd/mm/yy"
$ws2.Range("E1").Value = "Jeudi"
$ws2.Range("F1").Value = "Vendredi"

# ---- Row 3 : "Midi" banner ----------------------------------------------
foreach ($col in @("A", "B", "C", "D", "E", "F")) {
    $ws2.Range($col + "3").Value = "Midi"
}

# ---- Rows 4-8 : midi menu --------------------------------------------
$midiLabels = @("Entrée", "Plat chaud 1", "Plat chaud 2", "Dessert 1", "Dessert 2")
for ($i = 0; $i -lt $midiLabels.Length; $i++) {
    $r = 4 + $i
    $ws2.Range("A$r").Value = $midiLabels[$i]
    foreach ($col in @("B", "C", "D", "E", "F")) {
        $ws2.Range("$col$r").Value = $i + 1
    }
}
# Special menu item for the 04/04 midi Monday entrée
$ws2.Range("B4").Value = "Crêpes jambon"

# ---- Row 10 : "Soir" banner (only A-E, no Friday dinner service) --------
foreach ($col in @("A", "B", "C", "D", "E")) {
    $ws2.Range($col + "10").Value = "Soir"
}

# ---- Rows 11-15 : soir menu (only A-E) -----------------------------------
$soirLabels = @("Entrée", "Plat chaud 1", "Plat chaud 2", "Dessert 1", "Dessert 2")
for ($i = 0; $i -lt $soirLabels.Length; $i++) {
    $r = 11 + $i
    $ws2.Range("A$r").Value = $soirLabels[$i]
    foreach ($col in @("B", "C", "D", "E")) {
        $ws2.Range("$col$r").Value = $i + 6
    }
}

# ---- Row heights: the 04/04 sheet uses taller rows for the menu rows ----
foreach ($r in @(4, 5, 6, 7, 8, 11, 12, 13, 14, 15)) {
    $ws2.Rows.Item($r).RowHeight = 38.95
}
foreach ($r in @(1, 3, 10)) {
    $ws2.Rows.Item($r).RowHeight = 12.8
}

# ---- The original sheet's selection moves back to A1 (it is no longer ---
# ---- the tab that is selected when the workbook is reopened). -----------
$ws1.Range("A1").Select()

# ---- Select A1 on the new sheet and make it the active tab --------------
$ws2.Range("A1").Select()
$ws2.Activate()
